# "Finished c284130850 (not tested)"
#
# Row 40 of Sheet1 tracks task c284130850 (A40 = 284130850). Its status
# column moves from "i" (WIP/in-progress) to "w" (测试中 / in testing),
# since the work is done but has not been tested yet. The dependent
# summary formulas (WIP count, testing count, "other" subtotal) pick up
# the change automatically on recalculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B40").Value = "w"

# Leave the selection where the author ended up after making the edit.
$ws.Range("E35").Select()
